$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The daily scan re-ranked BTC-USD to the top row and refreshed all
# metrics with the latest (2025-12-06 21:20) run's numbers. RIOT, MARA and
# COIN each shift down a row; MSTR stays put at row 6.
#
# Column A holds the date as literal text in the source file (not a real
# date serial), so force a Text number format before assigning the string
# and reset the style back to Normal afterwards to avoid leaving a visible
# format change behind.
$ws.Range("A2:A6").NumberFormat = "@"

# Row 2: Bitcoin USD / BTC-USD (moved to top)
$ws.Range("A2").Value = "2025-12-06"
$ws.Range("B2").Value = "Bitcoin USD"
$ws.Range("C2").Value = "BTC-USD"
$ws.Range("D2").Value = 89623.46000000001
$ws.Range("E2").Value = 59.9
$ws.Range("F2").Value = 3.83
$ws.Range("G2").Value = 60
$ws.Range("H2").Value = 63
$ws.Range("I2").Value = 60
$ws.Range("J2").Value = 56
$ws.Range("K2").Value = 57.7
$ws.Range("L2").Value = "Pattern"
$ws.Range("M2").Value = "⛔ 관망하십시오."
$ws.Range("N2").Value = 52.28493729186943
$ws.Range("O2").Value = "⚪ 중립 구간"

# Row 3: Riot Platforms, Inc. / RIOT
$ws.Range("A3").Value = "2025-12-06"
$ws.Range("B3").Value = "Riot Platforms, Inc."
$ws.Range("C3").Value = "RIOT"
$ws.Range("D3").Value = 14.94
$ws.Range("E3").Value = 57.3
$ws.Range("F3").Value = -7.38
$ws.Range("G3").Value = 50
$ws.Range("H3").Value = 50
$ws.Range("I3").Value = 60
$ws.Range("J3").Value = 70
$ws.Range("K3").Value = 54.7
$ws.Range("L3").Value = "Pattern"
$ws.Range("M3").Value = "⛔ 관망하십시오."
$ws.Range("N3").Value = 52.28493729186943
$ws.Range("O3").Value = "⚪ 중립 구간"

# Row 4: MARA Holdings, Inc. / MARA
$ws.Range("A4").Value = "2025-12-06"
$ws.Range("B4").Value = "MARA Holdings, Inc."
$ws.Range("C4").Value = "MARA"
$ws.Range("D4").Value = 11.74
$ws.Range("E4").Value = 48.1
$ws.Range("F4").Value = -0.59
$ws.Range("G4").Value = 30
$ws.Range("H4").Value = 56
$ws.Range("I4").Value = 63
$ws.Range("J4").Value = 70
$ws.Range("K4").Value = 49.9
$ws.Range("L4").Value = "Pattern"
$ws.Range("M4").Value = "⛔ 관망하십시오."
$ws.Range("N4").Value = 52.28493729186943
$ws.Range("O4").Value = "⚪ 중립 구간"

# Row 5: Coinbase Global, Inc. / COIN
$ws.Range("A5").Value = "2025-12-06"
$ws.Range("B5").Value = "Coinbase Global, Inc."
$ws.Range("C5").Value = "COIN"
$ws.Range("D5").Value = 269.73
$ws.Range("E5").Value = 44.1
$ws.Range("F5").Value = -1.13
$ws.Range("G5").Value = 30
$ws.Range("H5").Value = 56
$ws.Range("I5").Value = 60
$ws.Range("J5").Value = 50
$ws.Range("K5").Value = 48.7
$ws.Range("L5").Value = "Pattern"
$ws.Range("M5").Value = "⛔ 관망하십시오."
$ws.Range("N5").Value = 52.28493729186943
$ws.Range("O5").Value = "⚪ 중립 구간"

# Row 6: Strategy Inc / MSTR (stays at the bottom)
$ws.Range("A6").Value = "2025-12-06"
$ws.Range("B6").Value = "Strategy Inc"
$ws.Range("C6").Value = "MSTR"
$ws.Range("D6").Value = 178.99
$ws.Range("E6").Value = 40.1
$ws.Range("F6").Value = 1.02
$ws.Range("G6").Value = 40
$ws.Range("H6").Value = 36
$ws.Range("I6").Value = 40
$ws.Range("J6").Value = 36
$ws.Range("K6").Value = 43.7
$ws.Range("L6").Value = "Pattern"
$ws.Range("M6").Value = "⛔ 관망하십시오."
$ws.Range("N6").Value = 52.28493729186943
$ws.Range("O6").Value = "⚪ 중립 구간"

# Restore the default (unstyled) look for column A now that the text
# values are safely stored as strings rather than parsed dates.
$ws.Range("A2:A6").Style = "Normal"
